$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the existing rows 10-16 to make room (in the underlying scheme
# ordering) for the three new spiral-scan rows that get inserted logically
# between "Ring Perpendicular to TD" and "NoRotation-tilt60deg". The sheet
# itself is just extended with new rows at the bottom (17-19) carrying the
# labels that used to sit in rows 10-16, and rows 10-16 are re-labelled to
# the new scheme order.
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# --- Add three new rows (17-19) for the HexGrid schemes that were bumped
# down by the insertion above, copying the formatting of the row above so
# the index column (A) keeps its bold/centered/bordered style.
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$rows = @(
    @{ Row = 17; Index = 15; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Index = 16; Label = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Index = 17; Label = "HexGrid-60degTilt5degRes" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.Index
    $ws.Cells.Item($rowNum, 2).Value = $r.Label
    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($rowNum, $c).Value = 1
    }
}

Write-Output "done"
